$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

# Header row (B2:E2)
Set-TextValue "B2" "Company"
Set-TextValue "C2" "Price"
Set-TextValue "D2" "Change"
Set-TextValue "E2" "Value  (Rs Cr.)"

# Row 3 - Adani Ports
Set-TextValue "B3" "Adani Ports"
Set-TextValue "C3" "904.20"
Set-TextValue "D3" "44.65"
Set-TextValue "E3" "1,380.40"

# Row 4 - ICICI Bank
Set-TextValue "B4" "ICICI Bank"
Set-TextValue "C4" "755.85"
Set-TextValue "D4" "3.65"
Set-TextValue "E4" "954.20"

# Row 5 - Reliance
Set-TextValue "B5" "Reliance"
Set-TextValue "C5" "2740.95"
Set-TextValue "D5" "45.95"
Set-TextValue "E5" "946.61"

# Row 6 - HDFC Bank
Set-TextValue "B6" "HDFC Bank"
Set-TextValue "C6" "1369.00"
Set-TextValue "D6" "3.25"
Set-TextValue "E6" "937.55"

# Row 7 - Axis Bank
Set-TextValue "B7" "Axis Bank"
Set-TextValue "C7" "780.40"
Set-TextValue "D7" "-0.30"
Set-TextValue "E7" "537.38"

# Rows 8 and 9 are removed entirely from the table
$ws.Range("A8:E9").Delete()
